$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.932.20"
$ws.Range("E2").Value = "  -1.40%  "
$ws.Range("D3").Value = "'1.888.42"
$ws.Range("E3").Value = "  -2.56%  "
$ws.Range("D4").Value = "'0.9999"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "'0.7343"
$ws.Range("E5").Value = "  -4.65%  "
$ws.Range("D6").Value = "'242.42"
$ws.Range("E6").Value = "  -1.33%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "'0.3108"
$ws.Range("E8").Value = "  -3.04%  "
$ws.Range("D9").Value = "'26.21"
$ws.Range("E9").Value = "  -5.64%  "
$ws.Range("E10").Value = "  -1.78%  "
$ws.Range("D11").Value = "'0.07937"
$ws.Range("E11").Value = "  -0.92%  "
$ws.Range("D12").Value = "'0.7661"
$ws.Range("E12").Value = "  -2.04%  "
$ws.Range("D13").Value = "'1.913.33"
$ws.Range("E13").Value = "  -1.07%  "
$ws.Range("D14").Value = "'5.229"
$ws.Range("E14").Value = "  -2.26%  "
$ws.Range("D15").Value = "'91.29"
$ws.Range("E15").Value = "  -3.54%  "
$ws.Range("E16").Value = "  -2.15%  "
$ws.Range("D17").Value = "'29.943.56"
$ws.Range("E17").Value = "  -1.38%  "
$ws.Range("D18").Value = "'5.753"
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("D19").Value = "'239.86"
$ws.Range("E19").Value = "  -6.12%  "
$ws.Range("D20").Value = "'0.000007754"
$ws.Range("E20").Value = "  -2.36%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("D22").Value = "'2.147.04"
$ws.Range("E22").Value = "  -1.93%  "
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("D24").Value = "'6.900"
$ws.Range("E24").Value = "  +2.71%  "
$ws.Range("D25").Value = "'9.296"
$ws.Range("E25").Value = "  -2.47%  "
$ws.Range("D26").Value = "'164.44"
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "'18.88"
$ws.Range("E27").Value = "  -1.00%  "
$ws.Range("D28").Value = "'0.1272"
$ws.Range("E28").Value = "  -5.15%  "
$ws.Range("D29").Value = "'2.017"
$ws.Range("E29").Value = "  -11.32%  "
$ws.Range("D30").Value = "'1.351"
$ws.Range("E30").Value = "  -1.02%  "
$ws.Range("E31").Value = "  +1.15%  "
$ws.Range("E32").Value = "  -2.45%  "
$ws.Range("D33").Value = "'4.084"
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("E34").Value = "  -1.32%  "
$ws.Range("D35").Value = "'1.278"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").Value = "'0.7366"
$ws.Range("E36").Value = "  -1.39%  "
$ws.Range("D37").Value = "'2.719"
$ws.Range("E37").Value = "  -2.19%  "
$ws.Range("D38").Value = "'0.01921"
$ws.Range("E38").Value = "  -1.70%  "
$ws.Range("E39").Value = "  -1.14%  "
$ws.Range("D40").Value = "'6.312"
$ws.Range("E40").Value = "  -1.72%  "
$ws.Range("D41").Value = "'74.59"
$ws.Range("E41").Value = "  -5.13%  "
$ws.Range("D42").Value = "'0.4450"
$ws.Range("E42").Value = "  -0.96%  "
$ws.Range("D43").Value = "'1.929"
$ws.Range("E43").Value = "  -2.19%  "
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").Value = "'0.8371"
$ws.Range("E45").Value = "  +0.48%  "
$ws.Range("D46").Value = "'7.590"
$ws.Range("E46").Value = "  +1.38%  "
$ws.Range("D47").Value = "'100.91"
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("D48").Value = "'9.776"
$ws.Range("D49").Value = "'37.12"
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("D50").Value = "'2.051.75"
$ws.Range("E50").Value = "  -1.77%  "
$ws.Range("D51").Value = "'942.85"
$ws.Range("E51").Value = "  -3.52%  "
